$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values pulled from the repulled dataset.
$updates = @{
    2  = -1
    3  = 1
    5  = -1
    6  = 1
    7  = -3
    12 = 0
    18 = -5
    19 = -3
    24 = 4
    25 = -2
    26 = -5
    27 = -3
    28 = -6
    29 = -4
    33 = -1
    39 = 5
    41 = -4
    42 = -2
    43 = 6
    44 = 4
    46 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
